$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two sampled rows ("RM 232" at row 26 and "SC 92" at row 28).
# Delete bottom-up so row indices for the earlier deletion stay valid.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# Re-roll which cells are "missing" (blank) vs populated for the remaining rows.
# Cells that became newly missing:
$ws.Range("F4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("D19").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("B26").ClearContents()
$ws.Range("C27").ClearContents()
$ws.Range("D27").ClearContents()
$ws.Range("B29").ClearContents()
$ws.Range("D29").ClearContents()
$ws.Range("F29").ClearContents()
$ws.Range("F30").ClearContents()
$ws.Range("F32").ClearContents()

# Cells that became newly populated (previously missing):
$ws.Range("C6").Value = 15.1
$ws.Range("D11").Value = -15.5
$ws.Range("F17").Value = 17.78
$ws.Range("C19").Value = 13.2
$ws.Range("C23").Value = 12.2
$ws.Range("D23").Value = -13.9
$ws.Range("F24").Value = 16.78
$ws.Range("D25").Value = -15.5
$ws.Range("B27").Value = -20.4
$ws.Range("F27").Value = 17
$ws.Range("F28").Value = 17.44
$ws.Range("C29").Value = 11.2
$ws.Range("D30").Value = -13.6
$ws.Range("D33").Value = -14.1
